# Change the user id in the "createUser" sheet from 1031 to 1032.
# The B2 and F2 cells contain formulas (CONCAT) that reference A2, so
# updating A2 automatically recalculates the dependent values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("createUser")
$ws.Range("A2").Value = 1032
$wb.Application.Calculate()
